$d = $word.ActiveDocument

# Replace paragraphs 48..48
$p1 = $d.Paragraphs.Item(48)
$p2 = $d.Paragraphs.Item(48)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="038658F7" w14:textId="11803CA4"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>boat.short</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# Replace paragraphs 46..46
$p1 = $d.Paragraphs.Item(46)
$p2 = $d.Paragraphs.Item(46)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="455C717B" w14:textId="5B3C2B7F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">{%p if watercraft </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>%}</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>WATERCRAFT</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# Replace paragraphs 43..43
$p1 = $d.Paragraphs.Item(43)
$p2 = $d.Paragraphs.Item(43)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2A5A9A43" w14:textId="4A907D44"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>ac.short</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# Replace paragraphs 41..41
$p1 = $d.Paragraphs.Item(41)
$p2 = $d.Paragraphs.Item(41)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="133F5D9C" w14:textId="41735DA0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">{%p if aircraft </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>%}</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>AIRCRAFT</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# Replace paragraphs 38..38
$p1 = $d.Paragraphs.Item(38)
$p2 = $d.Paragraphs.Item(38)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4EC76C05" w14:textId="22D8954F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>auto.short</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# Replace paragraphs 36..36
$p1 = $d.Paragraphs.Item(36)
$p2 = $d.Paragraphs.Item(36)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2CF504D7" w14:textId="7FA31C31"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{%p if automobile</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>%}AUTOMOBILES</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">{%p for </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>auto</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> in </w:t></w:r><w:r><w:t>automobile</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> %}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# Replace paragraphs 27..31
$p1 = $d.Paragraphs.Item(27)
$p2 = $d.Paragraphs.Item(31)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5FEE564A" w14:textId="77777777"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>acct.institution</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> }} {{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>acct.account_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}} x{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>acct.account_number</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="557136BD" w14:textId="7C97E707"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="4680"/></w:tabs></w:pPr><w:r><w:t>Participant/Owner:</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>acct.names</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_account</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="61F3229E" w14:textId="4FEB757C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="4680"/></w:tabs></w:pPr><w:r><w:t>Account type:</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>acct.account</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="2FF14F0D" w14:textId="60B0F749"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="4680"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">{%p if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>acct.is_employer_sponsored</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="4680"/></w:tabs></w:pPr><w:r><w:t>Employer:</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>acct.sponsor</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="4680"/></w:tabs></w:pPr><w:r><w:t>{%</w:t></w:r><w:r><w:t>p</w:t></w:r><w:r><w:t xml:space="preserve"> endif %}</w:t></w:r></w:p><w:p w14:paraId="734DD30A" w14:textId="77777777"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="decimal" w:pos="7920"/></w:tabs></w:pPr><w:r><w:t>Current balance ({{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>acct.balance</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_as_of_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}):</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>acct.current_balance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# Replace paragraphs 19..21
$p1 = $d.Paragraphs.Item(19)
$p2 = $d.Paragraphs.Item(21)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="504106BB" w14:textId="23D328A5"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>acct.institution</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"></w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}} {{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>acct.account_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">}} </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>x</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>acct.account_number</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="48DD1592" w14:textId="4A470784"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="4680"/></w:tabs></w:pPr><w:r><w:t>Account owner:</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>acct.names</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_account</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="37D399F3" w14:textId="5955C3D7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="decimal" w:pos="7920"/></w:tabs></w:pPr><w:r><w:t>Current balance</w:t></w:r><w:r><w:t xml:space="preserve"> ({{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>acct.balance</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_as_of_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}})</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>acct.current_balance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# Replace paragraphs 6..13
$p1 = $d.Paragraphs.Item(6)
$p2 = $d.Paragraphs.Item(13)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="424F83A0" w14:textId="3E08CDB5"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Address: </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>prop.address</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}}, {{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>prop.city</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}}, {{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>prop.county</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}}, {{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>prop.state</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}} {{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>prop.postal_code</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="53089440" w14:textId="6DE7A169"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="4680"/></w:tabs></w:pPr><w:r><w:t>Titled to:</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>prop.names</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_title</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="72248A03" w14:textId="214086C4"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="4680"/></w:tabs></w:pPr><w:r><w:t>Separate property claim</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>prop.sp</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_claimant</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="180BC267" w14:textId="7ADC57E3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="4680"/></w:tabs></w:pPr><w:r><w:t>Loan obligor:</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>prop.names</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_loan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="3BD76222" w14:textId="48A2AB0E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="decimal" w:pos="7920"/></w:tabs></w:pPr><w:r><w:t>Fair market value:</w:t></w:r><w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>{{ "</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>${:,.2f}".format(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prop.estimated_value</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)}}</w:t></w:r></w:p><w:p w14:paraId="7A52A5BD" w14:textId="72FEAF29"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="decimal" w:pos="7920"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Mortgage balance: </w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:r><w:t>"${:,.2f}</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>".format</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prop.lien_balances</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="7D98F8AB" w14:textId="2A572CB6"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="decimal" w:pos="7920"/></w:tabs></w:pPr><w:r><w:t>Property taxes due:</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:r><w:t>"${:,.2f}</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>".format</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prop.taxes_due</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="6D899504" w14:textId="2C08A032"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="decimal" w:pos="7920"/></w:tabs><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:t>Value of equity</w:t></w:r><w:r><w:tab/><w:t>{{</w:t></w:r><w:r><w:t>"${:,.2f}</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>".format</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prop.equity_value</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
